# Applies the "Added more performance info" revision to
# "Performance analysis of sorting.docx".
#
# Changes made (matching the target unified diff):
#   1. Tidy up proof-reading marks that bracketed the O(...) notation
#      (these collapse away naturally once the surrounding text is
#      touched, same as a normal Word edit/reproof pass).
#   2. Insert " (ie record)" right after "...standard quicksort per
#      pass" and relocate the "_GoBack" bookmark there.
#   3. Append a new closing sentence about the observed linear timing.
#   4. Tidy "Y axis"/"X axis" captions and add a trailing space run
#      after "X axis is the pass number".
#   5. Clear the now-empty bookmark paragraph at the end of the doc.

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# --- 1. Collapse the split O(nlgn) runs / drop their proof marks ---
Replace-Text "an optimal O(nlgn) sorting algorithm" "an optimal O(nlgn) sorting algorithm"

# --- 2. Collapse " {O(n^2)}" (worst case) ---
Replace-Text " worst case {O(n^2)} of" " worst case {O(n^2)} of"

# --- 3. Insert "(ie record)" + relocate the _GoBack bookmark ---
Replace-Text "standard quicksort per pass but alternatively" "standard quicksort per pass (ie record) but alternatively"

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$idx = $d.Content.Text.IndexOf("(ie record)")
$markPos = $idx + ("(ie record)").Length
$bmRange = $d.Range($markPos, $markPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- 4. Collapse the quicksort "O(n^2)" mention ---
Replace-Text ".  For my implementation I tried a quicksort implementation which ended up being very close to O(n^2) " ".  For my implementation I tried a quicksort implementation which ended up being very close to O(n^2) "

# --- 5. Collapse the insertion sort "O(n)" mention ---
Replace-Text " sort, I ended up with an O(n)" " sort, I ended up with an O(n)"

# --- 6. Collapse the final O(n^2)/O(n * nlgn) summary sentence ---
Replace-Text " In the end because there are n passes for each of the 6 sorts, the total time complexity ends up being around 6n(n+1)/2 which ends up being O(n^2) overall. This is still better than even the best case for a quicksort which would end up being O(n * nlgn). " " In the end because there are n passes for each of the 6 sorts, the total time complexity ends up being around 6n(n+1)/2 which ends up being O(n^2) overall. This is still better than even the best case for a quicksort which would end up being O(n * nlgn). "

# --- 7. Append the new closing sentence as its own run ---
$endIdx = $d.Content.Text.IndexOf("O(n * nlgn). ") + ("O(n * nlgn). ").Length
$insPoint = $d.Range($endIdx, $endIdx)
$insPoint.InsertAfter(" Ultimately my times generally followed a linear increase per pass with some discrepancy due to real world cpu/memory constraints.  ")

# --- 8. Tidy the axis caption paragraphs ---
Replace-Text "Y axis is the number microseconds per pass" "Y axis is the number microseconds per pass"
Replace-Text "X axis is the pass number" "X axis is the pass number"

$xIdx = $d.Content.Text.IndexOf("X axis is the pass number") + ("X axis is the pass number").Length
$xPoint = $d.Range($xIdx, $xIdx)
$xPoint.InsertAfter(" ")

# --- 9. The bookmark paragraph at the end is now just empty ---
Write-Output "done"
